$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 2).Value = 0.9565217391304348
$ws.Cells.Item(3, 3).Value = 44
$ws.Cells.Item(3, 4).Value = 44
$ws.Cells.Item(3, 8).Value = 2
$ws.Cells.Item(4, 2).Value = 0.8863636363636364
$ws.Cells.Item(4, 3).Value = 39
$ws.Cells.Item(4, 4).Value = 39
$ws.Cells.Item(4, 8).Value = 5
$ws.Cells.Item(4, 11).Value = 0.7846153846153846
$ws.Cells.Item(4, 12).Value = 51
$ws.Cells.Item(4, 13).Value = 51
$ws.Cells.Item(4, 17).Value = 14
$ws.Cells.Item(5, 1).Value = "however"
$ws.Cells.Item(5, 2).Value = 0.78125
$ws.Cells.Item(5, 3).Value = 50
$ws.Cells.Item(5, 4).Value = 50
$ws.Cells.Item(5, 8).Value = 14
$ws.Cells.Item(5, 11).Value = 0.7096774193548387
$ws.Cells.Item(5, 12).Value = 66
$ws.Cells.Item(5, 13).Value = 66
$ws.Cells.Item(5, 17).Value = 27
$ws.Cells.Item(6, 1).Value = "disappointed"
$ws.Cells.Item(6, 2).Value = 0.7580645161290323
$ws.Cells.Item(6, 3).Value = 141
$ws.Cells.Item(6, 4).Value = 141
$ws.Cells.Item(6, 8).Value = 45
$ws.Cells.Item(6, 11).Value = 0.5849056603773585
$ws.Cells.Item(6, 12).Value = 31
$ws.Cells.Item(6, 13).Value = 31
$ws.Cells.Item(6, 17).Value = 22
$ws.Cells.Item(7, 1).Value = "broke"
$ws.Cells.Item(7, 2).Value = 0.7087378640776699
$ws.Cells.Item(7, 3).Value = 146
$ws.Cells.Item(7, 4).Value = 146
$ws.Cells.Item(7, 8).Value = 60
$ws.Cells.Item(7, 10).Value = "thank"
$ws.Cells.Item(7, 11).Value = 0.4492753623188406
$ws.Cells.Item(7, 12).Value = 31
$ws.Cells.Item(7, 13).Value = 31
$ws.Cells.Item(7, 17).Value = 38
$ws.Cells.Item(8, 1).Value = "poor"
$ws.Cells.Item(8, 2).Value = 0.6901408450704225
$ws.Cells.Item(8, 3).Value = 49
$ws.Cells.Item(8, 4).Value = 49
$ws.Cells.Item(8, 8).Value = 22
$ws.Cells.Item(8, 10).Value = "excellent"
$ws.Cells.Item(8, 11).Value = 0.4375
$ws.Cells.Item(8, 17).Value = 36
$ws.Cells.Item(9, 2).Value = 0.6756756756756757
$ws.Cells.Item(9, 3).Value = 100
$ws.Cells.Item(9, 4).Value = 100
$ws.Cells.Item(9, 8).Value = 48
$ws.Cells.Item(9, 11).Value = 0.3680327868852459
$ws.Cells.Item(9, 12).Value = 449
$ws.Cells.Item(9, 13).Value = 449
$ws.Cells.Item(9, 17).Value = 771
$ws.Cells.Item(10, 1).Value = "junk"
$ws.Cells.Item(10, 2).Value = 0.6363636363636364
$ws.Cells.Item(10, 3).Value = 35
$ws.Cells.Item(10, 4).Value = 35
$ws.Cells.Item(10, 8).Value = 20
$ws.Cells.Item(10, 11).Value = 0.3285509325681492
$ws.Cells.Item(10, 12).Value = 229
$ws.Cells.Item(10, 13).Value = 229
$ws.Cells.Item(10, 17).Value = 468
$ws.Cells.Item(11, 1).Value = "water"
$ws.Cells.Item(11, 2).Value = 0.6190476190476191
$ws.Cells.Item(11, 3).Value = 26
$ws.Cells.Item(11, 4).Value = 26
$ws.Cells.Item(11, 8).Value = 16
$ws.Cells.Item(11, 11).Value = 0.3008298755186722
$ws.Cells.Item(11, 12).Value = 145
$ws.Cells.Item(11, 13).Value = 145
$ws.Cells.Item(11, 17).Value = 337
$ws.Cells.Item(12, 10).Value = "perfect"
$ws.Cells.Item(12, 11).Value = 0.2650602409638554
$ws.Cells.Item(12, 12).Value = 44
$ws.Cells.Item(12, 13).Value = 44
$ws.Cells.Item(12, 17).Value = 122
$ws.Cells.Item(13, 10).Value = "best"
$ws.Cells.Item(13, 11).Value = 0.25
$ws.Cells.Item(13, 12).Value = 30
$ws.Cells.Item(13, 13).Value = 30
$ws.Cells.Item(13, 17).Value = 90
$ws.Cells.Item(14, 1).Value = "probably"
$ws.Cells.Item(14, 2).Value = 0.5789473684210527
$ws.Cells.Item(14, 3).Value = 22
$ws.Cells.Item(14, 4).Value = 22
$ws.Cells.Item(14, 8).Value = 16
$ws.Cells.Item(15, 1).Value = "smaller"
$ws.Cells.Item(15, 2).Value = 0.5714285714285714
$ws.Cells.Item(15, 3).Value = 68
$ws.Cells.Item(15, 4).Value = 68
$ws.Cells.Item(15, 8).Value = 51
$ws.Cells.Item(16, 1).Value = "okay"
$ws.Cells.Item(16, 2).Value = 0.5185185185185185
$ws.Cells.Item(16, 3).Value = 28
$ws.Cells.Item(16, 4).Value = 28
$ws.Cells.Item(16, 8).Value = 26
$ws.Cells.Item(16, 10).Value = "fun"
$ws.Cells.Item(16, 11).Value = 0.112182296231376
$ws.Cells.Item(16, 12).Value = 128
$ws.Cells.Item(16, 13).Value = 128
$ws.Cells.Item(16, 17).Value = 1013
$ws.Cells.Item(17, 1).Value = "guess"
$ws.Cells.Item(17, 2).Value = 0.5185185185185185
$ws.Cells.Item(17, 3).Value = 28
$ws.Cells.Item(17, 4).Value = 28
$ws.Cells.Item(17, 8).Value = 26
$ws.Cells.Item(17, 10).Value = "christmas"
$ws.Cells.Item(17, 11).Value = 0.1044176706827309
$ws.Cells.Item(17, 12).Value = 26
$ws.Cells.Item(17, 13).Value = 26
$ws.Cells.Item(17, 14).Value = 1
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = $false
$ws.Cells.Item(17, 17).Value = 223
$ws.Cells.Item(18, 1).Value = "small"
$ws.Cells.Item(18, 2).Value = 0.5101449275362319
$ws.Cells.Item(18, 3).Value = 176
$ws.Cells.Item(18, 4).Value = 176
$ws.Cells.Item(18, 8).Value = 169
$ws.Cells.Item(18, 10).Value = "game"
$ws.Cells.Item(18, 11).Value = 0.05
$ws.Cells.Item(18, 12).Value = 77
$ws.Cells.Item(18, 13).Value = 78
$ws.Cells.Item(18, 14).Value = 0.99
$ws.Cells.Item(18, 15).Value = 0.01000000000000001
$ws.Cells.Item(18, 16).Value = $true
$ws.Cells.Item(18, 17).Value = 1463
$ws.Cells.Item(19, 1).Value = "plastic"
$ws.Cells.Item(19, 2).Value = 0.4803149606299212
$ws.Cells.Item(19, 3).Value = 61
$ws.Cells.Item(19, 4).Value = 61
$ws.Cells.Item(19, 8).Value = 66
$ws.Cells.Item(20, 1).Value = "broken"
$ws.Cells.Item(20, 2).Value = 0.4578313253012048
$ws.Cells.Item(20, 3).Value = 38
$ws.Cells.Item(20, 4).Value = 38
$ws.Cells.Item(20, 8).Value = 45
$ws.Cells.Item(21, 1).Value = "apart"
$ws.Cells.Item(21, 2).Value = 0.4315789473684211
$ws.Cells.Item(21, 3).Value = 41
$ws.Cells.Item(21, 4).Value = 41
$ws.Cells.Item(21, 8).Value = 54
$ws.Cells.Item(22, 1).Value = "difficult"
$ws.Cells.Item(22, 2).Value = 0.4044943820224719
$ws.Cells.Item(22, 3).Value = 36
$ws.Cells.Item(22, 4).Value = 36
$ws.Cells.Item(22, 8).Value = 53
$ws.Cells.Item(23, 1).Value = "paint"
$ws.Cells.Item(23, 2).Value = 0.3968253968253968
$ws.Cells.Item(23, 3).Value = 25
$ws.Cells.Item(23, 4).Value = 25
$ws.Cells.Item(23, 8).Value = 38
$ws.Cells.Item(24, 1).Value = "ok"
$ws.Cells.Item(24, 2).Value = 0.3359375
$ws.Cells.Item(24, 3).Value = 43
$ws.Cells.Item(24, 4).Value = 43
$ws.Cells.Item(24, 8).Value = 85
$ws.Cells.Item(25, 2).Value = 0.3267326732673267
$ws.Cells.Item(25, 3).Value = 66
$ws.Cells.Item(25, 4).Value = 66
$ws.Cells.Item(25, 8).Value = 136
$ws.Cells.Item(26, 1).Value = "cheap"
$ws.Cells.Item(26, 2).Value = 0.2843601895734597
$ws.Cells.Item(26, 3).Value = 60
$ws.Cells.Item(26, 4).Value = 60
$ws.Cells.Item(26, 8).Value = 151
$ws.Cells.Item(27, 2).Value = 0.2735042735042735
$ws.Cells.Item(27, 3).Value = 32
$ws.Cells.Item(27, 4).Value = 32
$ws.Cells.Item(27, 8).Value = 85
$ws.Cells.Item(28, 1).Value = "bit"
$ws.Cells.Item(28, 2).Value = 0.2551020408163265
$ws.Cells.Item(28, 3).Value = 25
$ws.Cells.Item(28, 4).Value = 25
$ws.Cells.Item(28, 8).Value = 73
$ws.Cells.Item(29, 1).Value = "size"
$ws.Cells.Item(29, 2).Value = 0.2422680412371134
$ws.Cells.Item(29, 3).Value = 47
$ws.Cells.Item(29, 4).Value = 47
$ws.Cells.Item(29, 8).Value = 147
$ws.Cells.Item(30, 1).Value = "money"
$ws.Cells.Item(30, 2).Value = 0.2310126582278481
$ws.Cells.Item(30, 3).Value = 73
$ws.Cells.Item(30, 4).Value = 73
$ws.Cells.Item(30, 8).Value = 243
$ws.Cells.Item(31, 2).Value = 0.2101449275362319
$ws.Cells.Item(31, 3).Value = 58
$ws.Cells.Item(31, 4).Value = 58
$ws.Cells.Item(31, 8).Value = 218
$ws.Cells.Item(32, 1).Value = "would"
$ws.Cells.Item(32, 2).Value = 0.1958456973293768
$ws.Cells.Item(32, 3).Value = 132
$ws.Cells.Item(32, 4).Value = 132
$ws.Cells.Item(32, 8).Value = 542
$ws.Cells.Item(33, 2).Value = 0.1867088607594937
$ws.Cells.Item(33, 3).Value = 59
$ws.Cells.Item(33, 4).Value = 59
$ws.Cells.Item(33, 8).Value = 257
$ws.Cells.Item(34, 1).Value = "back"
$ws.Cells.Item(34, 2).Value = 0.1857142857142857
$ws.Cells.Item(34, 3).Value = 26
$ws.Cells.Item(34, 4).Value = 26
$ws.Cells.Item(34, 8).Value = 114
$ws.Cells.Item(35, 1).Value = "could"
$ws.Cells.Item(35, 2).Value = 0.1719745222929936
$ws.Cells.Item(35, 3).Value = 27
$ws.Cells.Item(35, 4).Value = 27
$ws.Cells.Item(35, 8).Value = 130
$ws.Cells.Item(36, 1).Value = "price"
$ws.Cells.Item(36, 2).Value = 0.1609195402298851
$ws.Cells.Item(36, 3).Value = 56
$ws.Cells.Item(36, 4).Value = 56
$ws.Cells.Item(36, 8).Value = 292
$ws.Cells.Item(37, 2).Value = 0.1585903083700441
$ws.Cells.Item(37, 3).Value = 72
$ws.Cells.Item(37, 4).Value = 72
$ws.Cells.Item(37, 8).Value = 382
$ws.Cells.Item(38, 1).Value = "better"
$ws.Cells.Item(38, 2).Value = 0.1495327102803738
$ws.Cells.Item(38, 3).Value = 32
$ws.Cells.Item(38, 4).Value = 32
$ws.Cells.Item(38, 8).Value = 182
$ws.Cells.Item(39, 1).Value = "used"
$ws.Cells.Item(39, 2).Value = 0.1428571428571428
$ws.Cells.Item(39, 3).Value = 25
$ws.Cells.Item(39, 4).Value = 25
$ws.Cells.Item(39, 8).Value = 150
$ws.Cells.Item(40, 1).Value = "hard"
$ws.Cells.Item(40, 2).Value = 0.135
$ws.Cells.Item(40, 3).Value = 27
$ws.Cells.Item(40, 4).Value = 27
$ws.Cells.Item(40, 8).Value = 173
$ws.Cells.Item(41, 1).Value = "2"
$ws.Cells.Item(41, 2).Value = 0.1235955056179775
$ws.Cells.Item(41, 3).Value = 33
$ws.Cells.Item(41, 4).Value = 33
$ws.Cells.Item(41, 5).Value = 0
$ws.Cells.Item(41, 6).Value = 1
$ws.Cells.Item(41, 7).Value = $false
$ws.Cells.Item(41, 8).Value = 234
$ws.Cells.Item(42, 2).Value = 0.1214574898785425
$ws.Cells.Item(42, 3).Value = 30
$ws.Cells.Item(42, 4).Value = 31
$ws.Cells.Item(42, 5).Value = 0.03
$ws.Cells.Item(42, 6).Value = 0.97
$ws.Cells.Item(42, 7).Value = $true
$ws.Cells.Item(42, 8).Value = 217
$ws.Cells.Item(43, 1).Value = "use"
$ws.Cells.Item(43, 2).Value = 0.1178082191780822
$ws.Cells.Item(43, 3).Value = 43
$ws.Cells.Item(43, 4).Value = 43
$ws.Cells.Item(43, 8).Value = 322
$ws.Cells.Item(44, 1).Value = "expected"
$ws.Cells.Item(44, 2).Value = 0.1145833333333333
$ws.Cells.Item(44, 3).Value = 22
$ws.Cells.Item(44, 4).Value = 22
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 1
$ws.Cells.Item(44, 7).Value = $false
$ws.Cells.Item(44, 8).Value = 170
$ws.Cells.Item(45, 1).Value = "even"
$ws.Cells.Item(45, 2).Value = 0.1003584229390681
$ws.Cells.Item(45, 3).Value = 28
$ws.Cells.Item(45, 4).Value = 29
$ws.Cells.Item(45, 5).Value = 0.03
$ws.Cells.Item(45, 6).Value = 0.97
$ws.Cells.Item(45, 8).Value = 251
$ws.Cells.Item(46, 1).Value = "like"
$ws.Cells.Item(46, 2).Value = 0.07084019769357495
$ws.Cells.Item(46, 3).Value = 43
$ws.Cells.Item(46, 4).Value = 44
$ws.Cells.Item(46, 5).Value = 0.02
$ws.Cells.Item(46, 6).Value = 0.98
$ws.Cells.Item(46, 7).Value = $true
$ws.Cells.Item(46, 8).Value = 564
$ws.Cells.Item(47, 2).Value = 0.07042253521126761
$ws.Cells.Item(47, 3).Value = 25
$ws.Cells.Item(47, 4).Value = 25
$ws.Cells.Item(47, 8).Value = 330
$ws.Cells.Item(48, 1).Value = "little"
$ws.Cells.Item(48, 2).Value = 0.066815144766147
$ws.Cells.Item(48, 3).Value = 30
$ws.Cells.Item(48, 4).Value = 30
$ws.Cells.Item(48, 5).Value = 0
$ws.Cells.Item(48, 6).Value = 1
$ws.Cells.Item(48, 7).Value = $false
$ws.Cells.Item(48, 8).Value = 419
$ws.Cells.Item(49, 1).Value = "much"
$ws.Cells.Item(49, 2).Value = 0.05555555555555555
$ws.Cells.Item(49, 3).Value = 24
$ws.Cells.Item(49, 4).Value = 26
$ws.Cells.Item(49, 5).Value = 0.08
$ws.Cells.Item(49, 6).Value = 0.92
$ws.Cells.Item(49, 8).Value = 408
$ws.Cells.Item(50, 1).Value = "one"
$ws.Cells.Item(50, 2).Value = 0.04292929292929293
$ws.Cells.Item(50, 3).Value = 34
$ws.Cells.Item(50, 4).Value = 36
$ws.Cells.Item(50, 5).Value = 0.06
$ws.Cells.Item(50, 6).Value = 0.94
$ws.Cells.Item(50, 8).Value = 758

$ws.Range("J19:Q20").Clear()
